# Weekly update: insert a new pair of price records (Primera/Segunda) for
# "Femacal de La Calera - Piña" ahead of the existing history, pushing all
# the existing rows (308:406) down by two rows (-> 310:408).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 308:309 - everything currently at/after row
# 308 (through 406) shifts down to 310:408, dimension grows to A1:T408.
$ws.Rows("308:309").Insert()

# --- Row 308: new "Primera" record ---
$ws.Range("A308").Value = 3
$ws.Range("B308").Value = "Femacal de La Calera"
$ws.Range("C308").Value = "Coquimbo"
$ws.Range("D308").Value = 44524
$ws.Range("E308").Value = 5
$ws.Range("F308").Value = "Fruta"
$ws.Range("G308").Value = 100108
$ws.Range("H308").Value = "Tropicales y subtropicales"
$ws.Range("I308").Value = 100108005
$ws.Range("J308").Value = "Piña"
$ws.Range("K308").Value = "Caramelo"
$ws.Range("L308").Value = "Primera"
$ws.Range("M308").Value = 230
$ws.Range("N308").Value = 18000
$ws.Range("O308").Value = 18000
$ws.Range("P308").Value = 18000
$ws.Range("Q308").Value = "$/caja 12 unidades"
$ws.Range("R308").Value = "Ecuador"
$ws.Range("S308").Value = 1500
$ws.Range("T308").Value = 12

# --- Row 309: new "Segunda" record ---
$ws.Range("A309").Value = 3
$ws.Range("B309").Value = "Femacal de La Calera"
$ws.Range("C309").Value = "Coquimbo"
$ws.Range("D309").Value = 44524
$ws.Range("E309").Value = 5
$ws.Range("F309").Value = "Fruta"
$ws.Range("G309").Value = 100108
$ws.Range("H309").Value = "Tropicales y subtropicales"
$ws.Range("I309").Value = 100108005
$ws.Range("J309").Value = "Piña"
$ws.Range("K309").Value = "Caramelo"
$ws.Range("L309").Value = "Segunda"
$ws.Range("M309").Value = 60
$ws.Range("N309").Value = 18000
$ws.Range("O309").Value = 18000
$ws.Range("P309").Value = 18000
$ws.Range("Q309").Value = "$/caja 14 unidades"
$ws.Range("R309").Value = "Ecuador"
$ws.Range("S309").Value = 1286
$ws.Range("T309").Value = 14
